# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the content has
# moved from "In Translation" to "Ready for handoff", and bumps the
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
# to the new handoff-generation time.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$overview.Range("G2").Value = "2016-08-16 22:57:11" # Latest HO Xliff Generate Date

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"        # Status
$zhcn.Range("H2").Value = "2016-08-16 22:57:05"      # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"        # Status
$dede.Range("H2").Value = "2016-08-16 22:57:11"      # Latest Handoff Datetime

# --- Column widths: the longer "Ready for handoff" text widened the
# status-related columns (re-autofit by the report generator). The closest
# achievable width on this engine (quantized to 1/6 character) is used.
$overview.Columns.Item(5).ColumnWidth = 16.33  # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 16.33  # column F (de-de)
$zhcn.Columns.Item(3).ColumnWidth = 16.33       # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 16.33       # column C (Status)
